$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy cell formats to their new target styles (while row numbers are original) ---
$ws.Range("B4").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("A8").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B8").Copy()
$ws.Range("B4").PasteSpecial(-4122)

$ws.Range("D8").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Step 2: update text content ---
$ws.Range("A3").Value = "Firearm Purchase Prohibition"
$ws.Range("A4").Value = "Extension (code)"
$ws.Range("B4").Value = "Firearm Purchase Prohibition Code"
$ws.Range("D4").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/fppq-res-ext:FirearmPurchaseProhibition[@structures:id=../nc:ActivityPersonAssociation/nc:Activity/@structures:ref]/me-fpp-codes:FirearmPurchaseProhibitionCode"
$ws.Range("D10").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonBirthDate/nc:Date"
$ws.Range("D11").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonName/nc:PersonGivenName"
$ws.Range("D12").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonName/nc:PersonMiddleName"
$ws.Range("D13").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonName/nc:PersonSurName"
$ws.Range("D14").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/nc:PersonName/nc:PersonNameSuffixText"
$ws.Range("D15").Value = "/fppq-res-doc:FirearmPurchaseProhibitionQueryResults/fppq-res-ext:FirearmPurchaseProhibitionReport/nc:Person[@structures:id=../nc:ActivityPersonAssociation/nc:Person/@structures:ref]/j:PersonSexCode"

# --- Step 3: delete obsolete rows (Court Order block + old Extension row), descending order ---
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# --- Step 4: column width + selection ---
$ws.Columns.Item(1).ColumnWidth = 25.33
$ws.Range("C7").Select()
